# Updates the cryptos list: refreshed prices/volume percentages and
# fixes the Maker/RenderToken row order (rows 49-50).
# Values that look like plain numbers (e.g. "513.42") are written with a
# leading apostrophe so Excel keeps them as text, matching the source data
# which stores every Price/Volume cell as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.992.12'
$ws.Range('E2').Value = '  +7.10%  '
$ws.Range('D3').Value = '2.691.98'
$ws.Range('E3').Value = '  +11.84%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''513.42'
$ws.Range('E5').Value = '  +5.02%  '
$ws.Range('D6').Value = '''157.98'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').Value = '''0.604'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '2.688.85'
$ws.Range('E9').Value = '  +10.81%  '
$ws.Range('D10').Value = '''6.52'
$ws.Range('E10').Value = '  +3.26%  '
$ws.Range('E11').Value = '  +5.72%  '
$ws.Range('E12').Value = '  +3.95%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '3.135.59'
$ws.Range('E14').Value = '  +10.80%  '
$ws.Range('D15').Value = '60.982.31'
$ws.Range('E15').Value = '  +7.08%  '
$ws.Range('D16').Value = '''21.86'
$ws.Range('E16').Value = '  +5.52%  '
$ws.Range('D17').Value = '''0.0000141'
$ws.Range('E17').Value = '  +5.88%  '
$ws.Range('D18').Value = '2.688.91'
$ws.Range('E18').Value = '  +10.69%  '
$ws.Range('D19').Value = '''4.81'
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '''350.04'
$ws.Range('E20').Value = '  +8.04%  '
$ws.Range('D21').Value = '''10.54'
$ws.Range('E21').Value = '  +5.85%  '
$ws.Range('D22').Value = '''6.22'
$ws.Range('E22').Value = '  +3.87%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''60.59'
$ws.Range('E24').Value = '  +4.06%  '
$ws.Range('D25').Value = '''0.425'
$ws.Range('E25').Value = '  +4.69%  '
$ws.Range('D26').Value = '2.773.37'
$ws.Range('E26').Value = '  +9.67%  '
$ws.Range('E27').Value = '  +4.06%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '0.0₃0872'
$ws.Range('E29').Value = '  +11.84%  '
$ws.Range('E30').Value = '  +3.82%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').Value = '''19.64'
$ws.Range('E32').Value = '  +5.79%  '
$ws.Range('D33').Value = '''156.93'
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('D35').Value = '''5.75'
$ws.Range('E35').Value = '  +8.86%  '
$ws.Range('D36').Value = '''4.08'
$ws.Range('E36').Value = '  +9.62%  '
$ws.Range('E37').Value = '  +5.34%  '
$ws.Range('D38').Value = '''315.15'
$ws.Range('E38').Value = '  +17.59%  '
$ws.Range('E39').Value = '  +10.36%  '
$ws.Range('D40').Value = '''0.859'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('E41').Value = '  +6.87%  '
$ws.Range('D42').Value = '''0.838'
$ws.Range('E42').Value = '  +31.13%  '
$ws.Range('D43').Value = '''35.39'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('D44').Value = '''0.647'
$ws.Range('E44').Value = '  +8.58%  '
$ws.Range('D45').Value = '''0.0579'
$ws.Range('E45').Value = '  +8.83%  '
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').Value = '''20.12'
$ws.Range('E47').Value = '  +15.56%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.078.28'
$ws.Range('E49').Value = '  +10.80%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '''4.88'
$ws.Range('E50').Value = '  +6.57%  '
$ws.Range('E51').Value = '  +3.29%  '
